$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to stay text, matching the
# original inline-string storage (values like "228.86" would otherwise be
# auto-coerced to a number by Excel's smart-entry parsing).
$priceCells = "D2","D3","D5","D6","D7","D10","D11","D12","D13","D14","D15","D17","D18","D20","D22","D24","D25","D26","D27","D29","D31","D33","D34","D35","D38","D40","D41","D43","D44","D47","D48","D50"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "39.474.81"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3
$ws.Range("D3").Value = "2.166.19"
$ws.Range("E3").Value = "  +2.95%  "

# Row 5
$ws.Range("D5").Value = "228.86"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7
$ws.Range("D7").Value = "63.70"
$ws.Range("E7").Value = "  +2.14%  "

# Row 9
$ws.Range("E9").Value = "  +1.23%  "

# Row 10
$ws.Range("D10").Value = "0.0852"
$ws.Range("E10").Value = "  +1.17%  "

# Row 11
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("D12").Value = "16.06"
$ws.Range("E12").Value = "  +1.84%  "

# Row 13
$ws.Range("D13").Value = "2.487.70"
$ws.Range("E13").Value = "  +3.00%  "

# Row 14
$ws.Range("D14").Value = "22.05"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").Value = "0.815"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range("D17").Value = "2.166.63"
$ws.Range("E17").Value = "  +2.85%  "

# Row 18
$ws.Range("D18").Value = "39.501.70"
$ws.Range("E18").Value = "  +1.79%  "

# Row 19
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("D20").Value = "71.89"
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("E21").Value = "  +0.86%  "

# Row 22
$ws.Range("D22").Value = "229.68"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - becomes PancakeSwap (previously Toncoin)
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  +1.84%  "

# Row 25 - becomes Toncoin (previously PancakeSwap)
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").Value = "172.15"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "9.54"
$ws.Range("E27").Value = "  -1.21%  "

# Row 28
$ws.Range("E28").Value = "  +1.17%  "

# Row 29
$ws.Range("D29").Value = "19.87"
$ws.Range("E29").Value = "  +2.66%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  +5.64%  "

# Row 32
$ws.Range("E32").Value = "  +1.08%  "

# Row 33
$ws.Range("D33").Value = "4.62"
$ws.Range("E33").Value = "  +1.40%  "

# Row 34 - becomes THORChain (previously InternetComputer(DFINITY))
$ws.Range("B34").Value = "THORChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D34").Value = "7.08"
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - becomes InternetComputer(DFINITY) (previously THORChain)
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "4.73"
$ws.Range("E35").Value = "  -0.88%  "

# Row 36
$ws.Range("E36").Value = "  +0.40%  "

# Row 37
$ws.Range("E37").Value = "  +0.53%  "

# Row 38
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").Value = "  +0.26%  "

# Row 39
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").Value = "103.08"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +0.13%  "

# Row 42
$ws.Range("E42").Value = "  -1.54%  "

# Row 43
$ws.Range("D43").Value = "1.521.18"
$ws.Range("E43").Value = "  -0.69%  "

# Row 44
$ws.Range("D44").Value = "1.21"
$ws.Range("E44").Value = "  +1.69%  "

# Row 45
$ws.Range("E45").Value = "  +5.49%  "

# Row 46
$ws.Range("E46").Value = "  +0.63%  "

# Row 47
$ws.Range("D47").Value = "0.0925"
$ws.Range("E47").Value = "  +1.41%  "

# Row 48
$ws.Range("D48").Value = "4.27"
$ws.Range("E48").Value = "  +3.31%  "

# Row 49
$ws.Range("E49").Value = "  -1.01%  "

# Row 50
$ws.Range("D50").Value = "2.371.40"
$ws.Range("E50").Value = "  +3.01%  "

# Row 51
$ws.Range("E51").Value = "  -0.69%  "
